$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(41, 8).Value = 3000
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).Value = ""

$ws.Cells.Item(125, 8).Value = 1440.3334
$ws.Cells.Item(125, 9).Value = 1142.5
$ws.Cells.Item(125, 11).Value = 10282.5
$ws.Cells.Item(125, 13).Value = -7822.5

$ws.Cells.Item(132, 8).Value = 68525.125
$ws.Cells.Item(132, 9).Value = 3640.4
$ws.Cells.Item(132, 11).Value = 10921.2
$ws.Cells.Item(132, 13).Value = -8391.200000000001

$ws.Cells.Item(138, 8).Value = 2953.1777
$ws.Cells.Item(138, 9).Value = 1340
$ws.Cells.Item(138, 10).Value = 3681.7097
$ws.Cells.Item(138, 11).Value = 4020
$ws.Cells.Item(138, 12).Value = 11045.1291
$ws.Cells.Item(138, 13).Value = 1120
$ws.Cells.Item(138, 14).Value = -21325.1291

$ws.Cells.Item(141, 8).Value = 8303.200000000001
$ws.Cells.Item(141, 9).Value = 7559.1113
$ws.Cells.Item(141, 11).Value = 22677.3339
$ws.Cells.Item(141, 13).Value = -17497.3339

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4160.4116
$ws.Cells.Item(32, 9).Value = 4445.288
$ws.Cells.Item(32, 11).Value = 4445.288
$ws.Cells.Item(32, 13).Value = -4158.288

$ws.Cells.Item(45, 8).Value = 6998.4
$ws.Cells.Item(45, 9).Value = 997
$ws.Cells.Item(45, 11).Value = 997
$ws.Cells.Item(45, 13).Value = -620

$ws.Cells.Item(61, 8).Value = 18505.428
$ws.Cells.Item(61, 9).Value = 2557.6
$ws.Cells.Item(61, 11).Value = 2557.6
$ws.Cells.Item(61, 13).Value = -2345.6

$ws.Cells.Item(74, 8).Value = 3926.7932
$ws.Cells.Item(74, 9).Value = 3411.5417
$ws.Cells.Item(74, 11).Value = 3411.5417
$ws.Cells.Item(74, 13).Value = -2537.5417

$ws.Cells.Item(77, 8).Value = 3926.7932
$ws.Cells.Item(77, 9).Value = 3411.5417
$ws.Cells.Item(77, 11).Value = 17057.7085
$ws.Cells.Item(77, 13).Value = -12689.7085

$ws.Cells.Item(102, 8).Value = 37047308
$ws.Cells.Item(102, 10).Value = 83347810
$ws.Cells.Item(102, 12).Value = 83347810
$ws.Cells.Item(102, 14).Value = -83351054

$ws.Cells.Item(110, 8).Value = 4361
$ws.Cells.Item(110, 9).Value = 4408.8335
$ws.Cells.Item(110, 11).Value = 4408.8335
$ws.Cells.Item(110, 13).Value = -2363.8335

$ws.Cells.Item(119, 8).Value = 89998.5
$ws.Cells.Item(119, 10).Value = 89998.5
$ws.Cells.Item(119, 12).Value = 89998.5
$ws.Cells.Item(119, 14).Value = -99674.5

$ws.Cells.Item(122, 8).Value = 4460.2285
$ws.Cells.Item(122, 9).Value = 4220.8
$ws.Cells.Item(122, 10).Value = 5058.8
$ws.Cells.Item(122, 11).Value = 12662.4
$ws.Cells.Item(122, 12).Value = 15176.4
$ws.Cells.Item(122, 13).Value = -10212.4
$ws.Cells.Item(122, 14).Value = -20076.4

$ws.Cells.Item(136, 8).Value = 18505.428
$ws.Cells.Item(136, 9).Value = 2557.6
$ws.Cells.Item(136, 11).Value = 7672.799999999999
$ws.Cells.Item(136, 13).Value = -5122.799999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(40, 8).Value = 67995
$ws.Cells.Item(40, 10).Value = 67995
$ws.Cells.Item(40, 12).Value = 67995
$ws.Cells.Item(40, 14).Value = -68525

$ws.Cells.Item(94, 8).Value = 3006611.5
$ws.Cells.Item(94, 9).Value = 2986.2593
$ws.Cells.Item(94, 10).Value = 11116400
$ws.Cells.Item(94, 11).Value = 2986.2593
$ws.Cells.Item(94, 12).Value = 11116400
$ws.Cells.Item(94, 13).Value = -2535.2593
$ws.Cells.Item(94, 14).Value = -11117302

$ws.Cells.Item(99, 8).Value = 25055.176
$ws.Cells.Item(99, 9).Value = 29803.215
$ws.Cells.Item(99, 11).Value = 29803.215
$ws.Cells.Item(99, 13).Value = -28305.215

$ws.Cells.Item(107, 8).Value = 8422.087
$ws.Cells.Item(107, 9).Value = 9485.429
$ws.Cells.Item(107, 11).Value = 9485.429
$ws.Cells.Item(107, 13).Value = -7565.429

$ws.Cells.Item(138, 8).Value = 61051.535
$ws.Cells.Item(138, 10).Value = 61051.535
$ws.Cells.Item(138, 12).Value = 61051.535
$ws.Cells.Item(138, 14).Value = -71331.535

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 226.8125
$ws.Cells.Item(7, 10).Value = 280
$ws.Cells.Item(7, 12).Value = 280
$ws.Cells.Item(7, 14).Value = -506

$ws.Cells.Item(31, 8).Value = 3742.389
$ws.Cells.Item(31, 9).Value = 2709.111
$ws.Cells.Item(31, 11).Value = 2709.111
$ws.Cells.Item(31, 13).Value = -2414.111

$ws.Cells.Item(34, 8).Value = 3742.389
$ws.Cells.Item(34, 9).Value = 2709.111
$ws.Cells.Item(34, 11).Value = 2709.111
$ws.Cells.Item(34, 13).Value = -2507.111

$ws.Cells.Item(59, 8).Value = 665063
$ws.Cells.Item(59, 9).Value = 999999
$ws.Cells.Item(59, 11).Value = 999999
$ws.Cells.Item(59, 13).Value = -998854

$ws.Cells.Item(132, 8).Value = 2070.3914
$ws.Cells.Item(132, 9).Value = 1820.25
$ws.Cells.Item(132, 11).Value = 5460.75
$ws.Cells.Item(132, 13).Value = -2930.75

$ws.Cells.Item(134, 8).Value = 8394.429
$ws.Cells.Item(134, 9).Value = 8394.429
$ws.Cells.Item(134, 11).Value = 25183.287
$ws.Cells.Item(134, 13).Value = -22648.287

$ws.Cells.Item(140, 8).Value = 43749.875
$ws.Cells.Item(140, 10).Value = 43749.875
$ws.Cells.Item(140, 12).Value = 43749.875
$ws.Cells.Item(140, 14).Value = -54109.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(114, 8).Value = 2172
$ws.Cells.Item(114, 9).Value = 28
$ws.Cells.Item(114, 10).Value = 3244
$ws.Cells.Item(114, 11).Value = 84
$ws.Cells.Item(114, 12).Value = 9732
$ws.Cells.Item(114, 13).Value = 3170
$ws.Cells.Item(114, 14).Value = -16240

$ws.Cells.Item(117, 8).Value = 1109.5
$ws.Cells.Item(117, 9).Value = 425.66666
$ws.Cells.Item(117, 10).Value = 1793.3334
$ws.Cells.Item(117, 11).Value = 1276.99998
$ws.Cells.Item(117, 12).Value = 5380.0002
$ws.Cells.Item(117, 13).Value = 2165.00002
$ws.Cells.Item(117, 14).Value = -12264.0002

$ws.Cells.Item(121, 8).Value = 1710.6666
$ws.Cells.Item(121, 9).Value = 673.3333
$ws.Cells.Item(121, 10).Value = 2748
$ws.Cells.Item(121, 11).Value = 2019.9999
$ws.Cells.Item(121, 12).Value = 8244
$ws.Cells.Item(121, 13).Value = -709.9999
$ws.Cells.Item(121, 14).Value = -10864

$ws.Cells.Item(138, 8).Value = 7698.85
$ws.Cells.Item(138, 10).Value = 8237.846
$ws.Cells.Item(138, 12).Value = 24713.538
$ws.Cells.Item(138, 14).Value = -34993.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 7659.6943
$ws.Cells.Item(132, 10).Value = 9256
$ws.Cells.Item(132, 12).Value = 27768
$ws.Cells.Item(132, 14).Value = -32828

$ws.Cells.Item(137, 8).Value = 89890
$ws.Cells.Item(137, 10).Value = 89890
$ws.Cells.Item(137, 12).Value = 89890
$ws.Cells.Item(137, 14).Value = -100090

$ws.Cells.Item(140, 8).Value = 77070.89999999999
$ws.Cells.Item(140, 10).Value = 80000
$ws.Cells.Item(140, 12).Value = 80000
$ws.Cells.Item(140, 14).Value = -90360

$ws.Cells.Item(141, 8).Value = 70494.5
$ws.Cells.Item(141, 10).Value = 70494.5
$ws.Cells.Item(141, 12).Value = 70494.5
$ws.Cells.Item(141, 14).Value = -80854.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5561.5
$ws.Cells.Item(7, 9).Value = 5661.6
$ws.Cells.Item(7, 10).Value = 5311.25
$ws.Cells.Item(7, 11).Value = 5661.6
$ws.Cells.Item(7, 12).Value = 5311.25
$ws.Cells.Item(7, 13).Value = -5549.6
$ws.Cells.Item(7, 14).Value = -5535.25

$ws.Cells.Item(40, 8).Value = 7994
$ws.Cells.Item(40, 9).Value = 8399.6
$ws.Cells.Item(40, 11).Value = 8399.6
$ws.Cells.Item(40, 13).Value = -8263.6

$ws.Cells.Item(46, 8).Value = 3344.3635
$ws.Cells.Item(46, 10).Value = 3528.8
$ws.Cells.Item(46, 12).Value = 3528.8
$ws.Cells.Item(46, 14).Value = -3904.8

$ws.Cells.Item(61, 8).Value = 2071.3333
$ws.Cells.Item(61, 9).Value = 2032.3636
$ws.Cells.Item(61, 10).Value = 2500
$ws.Cells.Item(61, 11).Value = 2032.3636
$ws.Cells.Item(61, 12).Value = 2500
$ws.Cells.Item(61, 13).Value = -1830.3636
$ws.Cells.Item(61, 14).Value = -2904

$ws.Cells.Item(103, 8).Value = 45840.79
$ws.Cells.Item(103, 10).Value = 46443.055
$ws.Cells.Item(103, 12).Value = 46443.055
$ws.Cells.Item(103, 14).Value = -48787.055

$ws.Cells.Item(113, 8).Value = 2071.3333
$ws.Cells.Item(113, 9).Value = 2032.3636
$ws.Cells.Item(113, 10).Value = 2500
$ws.Cells.Item(113, 11).Value = 2032.3636
$ws.Cells.Item(113, 12).Value = 2500
$ws.Cells.Item(113, 13).Value = 137.6364000000001
$ws.Cells.Item(113, 14).Value = -6840

$ws.Cells.Item(126, 8).Value = 5561.5
$ws.Cells.Item(126, 9).Value = 5661.6
$ws.Cells.Item(126, 10).Value = 5311.25
$ws.Cells.Item(126, 11).Value = 16984.8
$ws.Cells.Item(126, 12).Value = 15933.75
$ws.Cells.Item(126, 13).Value = -14514.8
$ws.Cells.Item(126, 14).Value = -20873.75

$ws.Cells.Item(136, 8).Value = 4401.3
$ws.Cells.Item(136, 10).Value = 10000
$ws.Cells.Item(136, 12).Value = 30000
$ws.Cells.Item(136, 14).Value = -35100

$ws.Cells.Item(139, 8).Value = 54998
$ws.Cells.Item(139, 10).Value = 40000
$ws.Cells.Item(139, 12).Value = 40000
$ws.Cells.Item(139, 14).Value = -50280

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = 0
$ws.Cells.Item(5, 13).Value = ""
$ws.Cells.Item(5, 14).Value = ""

$ws.Cells.Item(122, 8).Value = 1577.25
$ws.Cells.Item(122, 10).Value = 1587.25
$ws.Cells.Item(122, 12).Value = 4761.75
$ws.Cells.Item(122, 14).Value = -9661.75

$ws.Cells.Item(126, 8).Value = 2777.6667
$ws.Cells.Item(126, 9).Value = 2599.8333
$ws.Cells.Item(126, 10).Value = 3133.3333
$ws.Cells.Item(126, 11).Value = 7799.499899999999
$ws.Cells.Item(126, 12).Value = 9399.999899999999
$ws.Cells.Item(126, 13).Value = -5329.499899999999
$ws.Cells.Item(126, 14).Value = -14339.9999

$ws.Cells.Item(136, 8).Value = 5007.357
$ws.Cells.Item(136, 9).Value = 5009
$ws.Cells.Item(136, 11).Value = 15027
$ws.Cells.Item(136, 13).Value = -12477

$ws.Cells.Item(140, 8).Value = 79999
$ws.Cells.Item(140, 10).Value = 79999
$ws.Cells.Item(140, 12).Value = 79999
$ws.Cells.Item(140, 14).Value = -90359

$ws.Cells.Item(141, 8).Value = 82630.10000000001
$ws.Cells.Item(141, 10).Value = 82630.10000000001
$ws.Cells.Item(141, 12).Value = 82630.10000000001
$ws.Cells.Item(141, 14).Value = -92990.10000000001
